$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Insert "Gabon" as a new shared string entry right between
#        "Trinidad yTobago" and "Etiopia" (i.e. shift the text shown by
#        the existing row that used to read "Etiopia" and friends).
#        Because the cells in column A reference shared strings by index,
#        and Excel's COM model only exposes cell text (not raw shared
#        string indices), we simply set the text of the affected cells so
#        that the saved workbook's shared string table ends up reordered
#        the same way: "Gabon" inserted right after "Trinidad yTobago".

$ws.Range("A138").Value = "Gabon"
$ws.Range("A139").Value = "Etiopia"
$ws.Range("A140").Value = "Aruba"
$ws.Range("A141").Value = "Guayana Francesa"

# --- 2. Update the statistics cells.

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 698851
$ws.Range("C4").Value = 21281
$ws.Range("D4").Value = 59328
$ws.Range("E4").Value = 602681
$ws.Range("G4").Value = 2225
$ws.Range("H4").Value = 36842

# Row 8 - Alemania
$ws.Range("B8").Value = 139702
$ws.Range("C8").Value = 2004
$ws.Range("E8").Value = 53699

# Row 18 - Suiza
$ws.Range("E18").Value = 9351
$ws.Range("G18").Value = 46
$ws.Range("H18").Value = 1327

# Row 48 - Republica Dominicana
$ws.Range("D48").Value = 268
$ws.Range("E48").Value = 3658

# Row 131 - Ruanda
$ws.Range("B131").Value = 143
$ws.Range("C131").Value = 5
$ws.Range("D131").Value = 65

# Row 138 - now Gabon (new data)
$ws.Range("B138").Value = 108
$ws.Range("C138").Value = 28
$ws.Range("D138").Value = 7
$ws.Range("E138").Value = 100
$ws.Range("F138").Value = 0
$ws.Range("H138").Value = 1

# Row 139 - now Etiopia (old Etiopia data)
$ws.Range("B139").Value = 96
$ws.Range("C139").Value = 4
$ws.Range("D139").Value = 15
$ws.Range("E139").Value = 78
$ws.Range("F139").Value = 1
$ws.Range("H139").Value = 3

# Row 140 - now Aruba (old Aruba data)
$ws.Range("B140").Value = 96
$ws.Range("C140").Value = 1
$ws.Range("D140").Value = 43
$ws.Range("E140").Value = 51
$ws.Range("F140").Value = 1
$ws.Range("H140").Value = 2

# Row 141 - now Guayana Francesa (old Guayana Francesa data)
$ws.Range("B141").Value = 96
$ws.Range("C141").Value = 10
$ws.Range("D141").Value = 61
$ws.Range("E141").Value = 35
$ws.Range("F141").Value = 2
$ws.Range("H141").Value = 0

# Row 156 - Bahamas
$ws.Range("B156").Value = 54
$ws.Range("C156").Value = 1
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 9
